$d = $word.ActiveDocument

# --- Paragraph 4: ">>>  your stuff after this line >>>" ---
# Originally split across three runs ("...", "...", "...") with
# proofErr gramStart/gramEnd markers in between. Normalize it down to
# a single clean run by rewriting the paragraph's text. Because the
# concatenated text is already identical to the desired text, we first
# stage a temporary value (forcing the run-rewrite) and then set the
# final text, which collapses everything into one run and drops the
# proofErr markers.
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.Text = "__TMP__"
$p4b = $d.Paragraphs(4)
$r4b = $p4b.Range
$r4b.End = $r4b.End - 1
$r4b.Text = ">>>  your stuff after this line >>>"

# --- Paragraph 5: "Baz changes" (with _GoBack bookmark) ---
# Replace with the new "version management" commentary paragraph.
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$r5.End = $r5.End - 1
$r5.Text = "The version management system is quite a good collaboration tool that maintains the consistency of software developed by agroup of developer spread over a large geographical ares, thanks to Linus torvald for this great innovation"

# --- New paragraph: "Changes by " ---
$p5c = $d.Paragraphs(5)
$p5c.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$r6.End = $r6.End - 1
$r6.Text = "Changes by "

# --- Two additional trailing empty paragraphs (2 -> 4 empty paragraphs) ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$lastPara2 = $d.Paragraphs($d.Paragraphs.Count)
$lastPara2.Range.InsertParagraphAfter()
